$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row (row 78) appends the latest resale-number scrape for 2024-01-19 22:53:29.
# Columns A-D are textual (date/time/weekday/week-code) and must stay as literal
# text rather than being auto-converted to Excel date/time/number serials, so we
# force the range to Text format before writing, then clear the formatting again
# afterwards so the cells end up with no explicit style (matching the rest of the
# sheet's plain data rows).
$textRange = $ws.Range("A78:D78")
$textRange.NumberFormat = "@"
$ws.Range("A78").Value = "2024-01-19"
$ws.Range("B78").Value = "22:53:29"
$ws.Range("C78").Value = "Friday"
$ws.Range("D78").Value = "02"
$textRange.ClearFormats()

# Columns E-T are the per-city numeric resale counts (plain numbers, -1 = no data).
$ws.Range("E78").Value = 138256
$ws.Range("F78").Value = 140586
$ws.Range("G78").Value = 171568
$ws.Range("H78").Value = 148804
$ws.Range("I78").Value = -1
$ws.Range("J78").Value = 122580
$ws.Range("K78").Value = 223606
$ws.Range("L78").Value = 255302
$ws.Range("M78").Value = 185371
$ws.Range("N78").Value = 110324
$ws.Range("O78").Value = 41405
$ws.Range("P78").Value = 30922
$ws.Range("Q78").Value = 73600
$ws.Range("R78").Value = -1
$ws.Range("S78").Value = 42892
$ws.Range("T78").Value = -1
